$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing rows 2-4 down to 3-5.
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the Romanian Liga II match data.
# The Date column ("2025-10-07") needs a leading apostrophe so Excel
# stores it as literal text instead of auto-converting it to a date
# serial number; ClearFormats() afterwards strips the resulting
# quote-prefix style (and any style inherited from the Insert) so the
# row ends up with the same "no explicit style" cells as the rest of
# the data rows.
$ws.Cells.Item(2, 1).Value = "Romanian Liga II"
$ws.Cells.Item(2, 2).Value = "'2025-10-07"
$ws.Cells.Item(2, 3).Value = "14:00:00"
$ws.Cells.Item(2, 4).Value = "ASA Targu Mures"
$ws.Cells.Item(2, 5).Value = "Gloria Bistrita-Nasaud"
$ws.Cells.Item(2, 6).Value = 1.04
$ws.Cells.Item(2, 7).Value = 1000
$ws.Cells.Item(2, 8).Value = 1.04
$ws.Cells.Item(2, 9).Value = 1000
$ws.Cells.Item(2, 10).Value = 1.01
$ws.Cells.Item(2, 11).Value = 1000
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = 0
$ws.Cells.Item(2, 14).Value = 0
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 16).Value = 1.24
$ws.Cells.Item(2, 17).Value = 1.01
$ws.Cells.Item(2, 18).Value = 0
$ws.Cells.Item(2, 19).Value = 0
$ws.Cells.Item(2, 20).Value = 0
$ws.Cells.Item(2, 21).Value = 0
$ws.Cells.Item(2, 22).Value = 0
$ws.Cells.Item(2, 23).Value = 0
$ws.Cells.Item(2, 24).Value = 0
$ws.Cells.Item(2, 25).Value = 0
$ws.Cells.Item(2, 26).Value = 0
$ws.Cells.Item(2, 27).Value = 0
$ws.Cells.Item(2, 28).Value = 0
$ws.Cells.Item(2, 29).Value = 0
$ws.Cells.Item(2, 30).Value = 0
$ws.Cells.Item(2, 31).Value = 0
$ws.Cells.Item(2, 32).Value = 0
$ws.Cells.Item(2, 33).Value = 0
$ws.Cells.Item(2, 34).Value = 0
$ws.Cells.Item(2, 35).Value = 0
$ws.Cells.Item(2, 36).Value = 0
$ws.Cells.Item(2, 37).Value = 0
$ws.Cells.Item(2, 38).Value = 0
$ws.Cells.Item(2, 39).Value = 0
$ws.Cells.Item(2, 40).Value = 0
$ws.Cells.Item(2, 41).Value = 0

$ws.Rows.Item(2).ClearFormats()
